# Weekly update: Fruta / hortaliza, semanal
# Applies the data refresh for "Fruta, Terminal La Palmera de La Serena - Membrillo":
#  - Rows 7-25 get updated Fecha/Volumen/Precio/Origen values (new week's data shifted in)
#  - Two additional observations are appended as new rows 26-27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = 44659
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = 295000
$ws.Range("O7").Value = 300000
$ws.Range("P7").Value = 297500
$ws.Range("S7").Value = 661
$ws.Range("D8").Value = 44659
$ws.Range("L8").Value = 'Segunda'
$ws.Range("D9").Value = 44627
$ws.Range("M9").Value = 16
$ws.Range("N9").Value = 405000
$ws.Range("O9").Value = 410000
$ws.Range("P9").Value = 407500
$ws.Range("S9").Value = 906
$ws.Range("D10").Value = 44273
$ws.Range("M10").Value = 10
$ws.Range("D11").Value = 44273
$ws.Range("D12").Value = 44316
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 255000
$ws.Range("O12").Value = 260000
$ws.Range("P12").Value = 257500
$ws.Range("S12").Value = 572
$ws.Range("D13").Value = 44316
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 20
$ws.Range("N13").Value = 225000
$ws.Range("O13").Value = 230000
$ws.Range("P13").Value = 227500
$ws.Range("S13").Value = 506
$ws.Range("D14").Value = 44658
$ws.Range("N14").Value = 295000
$ws.Range("O14").Value = 300000
$ws.Range("P14").Value = 297500
$ws.Range("S14").Value = 661
$ws.Range("D15").Value = 44658
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 16
$ws.Range("N15").Value = 255000
$ws.Range("O15").Value = 260000
$ws.Range("P15").Value = 257500
$ws.Range("S15").Value = 572
$ws.Range("D16").Value = 44630
$ws.Range("M16").Value = 16
$ws.Range("N16").Value = 400000
$ws.Range("O16").Value = 405000
$ws.Range("P16").Value = 402500
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 894
$ws.Range("D17").Value = 44649
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 18
$ws.Range("N17").Value = 330000
$ws.Range("O17").Value = 340000
$ws.Range("P17").Value = 335000
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 744
$ws.Range("D18").Value = 44295
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 20
$ws.Range("N18").Value = 225000
$ws.Range("O18").Value = 230000
$ws.Range("P18").Value = 227500
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 506
$ws.Range("D19").Value = 44295
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 16
$ws.Range("N19").Value = 195000
$ws.Range("O19").Value = 200000
$ws.Range("P19").Value = 197500
$ws.Range("R19").Value = 'Región Metropolitana'
$ws.Range("S19").Value = 439
$ws.Range("D20").Value = 44648
$ws.Range("L20").Value = 'Especial'
$ws.Range("M20").Value = 16
$ws.Range("N20").Value = 335000
$ws.Range("O20").Value = 340000
$ws.Range("P20").Value = 337500
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 750
$ws.Range("D21").Value = 44628
$ws.Range("M21").Value = 14
$ws.Range("N21").Value = 400000
$ws.Range("O21").Value = 410000
$ws.Range("P21").Value = 405000
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 900
$ws.Range("D22").Value = 44635
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = 300000
$ws.Range("O22").Value = 310000
$ws.Range("P22").Value = 305000
$ws.Range("R22").Value = 'Región Metropolitana'
$ws.Range("S22").Value = 678
$ws.Range("D23").Value = 44634
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 10
$ws.Range("R23").Value = 'Región Metropolitana'
$ws.Range("D24").Value = 44622
$ws.Range("M24").Value = 16
$ws.Range("N24").Value = 410000
$ws.Range("O24").Value = 420000
$ws.Range("P24").Value = 415000
$ws.Range("R24").Value = 'Región de O''Higgins'
$ws.Range("S24").Value = 922
$ws.Range("L25").Value = 'Especial'
$ws.Range("N25").Value = 305000
$ws.Range("O25").Value = 310000
$ws.Range("P25").Value = 307500
$ws.Range("S25").Value = 683

$ws.Range("A26").Value = 8
$ws.Range("B26").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C26").Value = 'Coquimbo'
$ws.Range("D26").Value = 44309
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = 'Fruta'
$ws.Range("G26").Value = 100104
$ws.Range("H26").Value = 'Frutos de pepita'
$ws.Range("I26").Value = 100104003
$ws.Range("J26").Value = 'Membrillo'
$ws.Range("K26").Value = 'Champion'
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 20
$ws.Range("N26").Value = 285000
$ws.Range("O26").Value = 290000
$ws.Range("P26").Value = 287500
$ws.Range("Q26").Value = '$/bins (450 kilos)'
$ws.Range("R26").Value = 'Provincia de Cachapoal'
$ws.Range("S26").Value = 639
$ws.Range("T26").Value = 450
$ws.Range("A27").Value = 8
$ws.Range("B27").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C27").Value = 'Coquimbo'
$ws.Range("D27").Value = 44309
$ws.Range("E27").Value = 4
$ws.Range("F27").Value = 'Fruta'
$ws.Range("G27").Value = 100104
$ws.Range("H27").Value = 'Frutos de pepita'
$ws.Range("I27").Value = 100104003
$ws.Range("J27").Value = 'Membrillo'
$ws.Range("K27").Value = 'Champion'
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 20
$ws.Range("N27").Value = 255000
$ws.Range("O27").Value = 260000
$ws.Range("P27").Value = 257500
$ws.Range("Q27").Value = '$/bins (450 kilos)'
$ws.Range("R27").Value = 'Provincia de Cachapoal'
$ws.Range("S27").Value = 572
$ws.Range("T27").Value = 450

# Match the date formatting used by the rest of the "Fecha" column (style s="2")
# for the two newly appended rows.
$ws.Range("D26").NumberFormat = $ws.Range("D25").NumberFormat
$ws.Range("D27").NumberFormat = $ws.Range("D25").NumberFormat
